# Apriori / snopes crawler update:
#  - add a new "Sheet2" (placed after Sheet1, becomes the active sheet)
#    which will hold the summarized ai-pattern / fact-check crawl output.
#  - write its header row: id, A, C, num_Bs, exist, non-exist, meaningful, unmeaningful

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so the tab order matches Sheet1, Sheet2
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$headers = @("id", "A", "C", "num_Bs", "exist", "non-exist", "meaningful", "unmeaningful")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $ws2.Cells.Item(1, $col + 1).Value = $headers[$col]
}

# Match the author's last on-screen selection for the new sheet
$ws2.Range("I5").Select()
